$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "IOYhp_0105719"
$ws.Range("D2").Value = "kjqKUGHAUp"
$ws.Range("N2").Value = "UCN 10514"
$ws.Range("K2").Value = "VhnzziLRir"
$ws.Range("M2").Value = "VhnzziLRir"
